# Apply "Added System test cases till req-36" edit
# (cell writes are ordered to reproduce the shared-string insertion order of the
#  original authored workbook: F37,H37,I37/J37, F38,J38,I38, F40, F39,G39,H39,I39/J39)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases & Results")

# --- Row 37 (TestCase 35 / REQ-34) ---
$ws.Cells.Item(37,5).Value  = "High Impact"
$ws.Cells.Item(37,6).Value  = "Test if If No match is found in REQ-31 App shall show a notification to user for wrong password"
$ws.Cells.Item(37,7).Value  = "App must be installed"
$ws.Cells.Item(37,8).Value  = "Open the app and input random characters into the adm and password fields"
$ws.Cells.Item(37,9).Value  = "App shows a snackbar indicating wrong login details "
$ws.Cells.Item(37,10).Value = "App shows a snackbar indicating wrong login details "
$ws.Rows.Item(37).RowHeight = 57.6

# --- Row 38 (TestCase 36 / REQ-35) ---
$ws.Cells.Item(38,5).Value  = "High Impact"
$ws.Cells.Item(38,6).Value  = "Test that From req-33, user should be prompted by 2 buttons 1 to reserve book and 1 to view borrowed books "
$ws.Cells.Item(38,7).Value  = "App must be installed and account is created"
$ws.Cells.Item(38,8).Value  = "Open the app and input the adm field ""P2426082"" and password ""123"""
$ws.Cells.Item(38,10).Value = "App must show 2 buttons in homepage"
$ws.Cells.Item(38,9).Value  = "App must show 2 buttons in homepage one to reserve and 1 to view loaned books"
$ws.Rows.Item(38).RowHeight = 57.6

# --- Row 40 (TestCase 38 / REQ-37): only Description/Test Summary filled in so far ---
$ws.Cells.Item(40,6).Value = "Test that"

# --- Row 39 (TestCase 37 / REQ-36) ---
$ws.Cells.Item(39,5).Value  = "Mid Impact"
$ws.Cells.Item(39,6).Value  = "Test that If user selects to view borrowed books, A list of reserved books and their loaned dates, title and location should be shown to user"
$ws.Cells.Item(39,7).Value  = "App must be installed and account is created and logged in"
$ws.Cells.Item(39,8).Value  = "presses the view the view loaned books button"
$ws.Cells.Item(39,9).Value  = "App must show a list of boxes with the loaned/reserved book's date on loan, title and location"
$ws.Cells.Item(39,10).Value = "App must show a list of boxes with the loaned/reserved book's date on loan, title and location"
$ws.Rows.Item(39).RowHeight = 72

# --- Update the visible selection / scroll position to match the new active cell ---
$ws.Activate()
$ws.Range("J39").Select()
$excel.ActiveWindow.ScrollRow = 37
